$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "30.830.52"
Set-TextValue "E2" "  -0.44%  "
Set-TextValue "D3" "1.937.87"
Set-TextValue "E3" "  -0.67%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "243.47"
Set-TextValue "E5" "  -0.63%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.16%  "
Set-TextValue "D7" "0.4885"
Set-TextValue "E7" "  +0.23%  "
Set-TextValue "D8" "0.2948"
Set-TextValue "E8" "  -0.46%  "
Set-TextValue "D9" "0.06893"
Set-TextValue "E9" "  +0.84%  "
Set-TextValue "D10" "19.32"
Set-TextValue "E10" "  +1.11%  "
Set-TextValue "D11" "104.78"
Set-TextValue "E11" "  -2.25%  "
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.939.98"
Set-TextValue "E12" "  -0.47%  "
Set-TextValue "B13" "TRON"
Set-TextValue "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.07780"
Set-TextValue "E13" "  +0.62%  "
Set-TextValue "D14" "5.343"
Set-TextValue "E14" "  -2.24%  "
Set-TextValue "D15" "0.7012"
Set-TextValue "E15" "  -0.63%  "
Set-TextValue "D16" "273.32"
Set-TextValue "E16" "  -2.92%  "
Set-TextValue "D17" "30.814.66"
Set-TextValue "E17" "  -0.59%  "
Set-TextValue "D18" "0.000007734"
Set-TextValue "E18" "  -0.10%  "
Set-TextValue "B19" "Uniswap"
Set-TextValue "C19" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D19" "5.639"
Set-TextValue "E19" "  +1.91%  "
Set-TextValue "B20" "Avalanche"
Set-TextValue "C20" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D20" "13.10"
Set-TextValue "E20" "  -1.17%  "
Set-TextValue "E21" "  +0.17%  "
Set-TextValue "D22" "2.182.60"
Set-TextValue "E22" "  -1.40%  "
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "6.533"
Set-TextValue "E24" "  +0.51%  "
Set-TextValue "D25" "9.828"
Set-TextValue "E25" "  +0.09%  "
Set-TextValue "D26" "165.37"
Set-TextValue "E26" "  -2.25%  "
Set-TextValue "D27" "19.65"
Set-TextValue "E27" "  -1.54%  "
Set-TextValue "D28" "2.160"
Set-TextValue "E28" "  -2.71%  "
Set-TextValue "D29" "0.1040"
Set-TextValue "E29" "  -1.29%  "
Set-TextValue "E30" "  -2.04%  "
Set-TextValue "D31" "4.631"
Set-TextValue "E31" "  +1.42%  "
Set-TextValue "D32" "1.558"
Set-TextValue "E32" "  -1.55%  "
Set-TextValue "D33" "4.389"
Set-TextValue "E33" "  -1.83%  "
Set-TextValue "D34" "0.04895"
Set-TextValue "E34" "  -1.55%  "
Set-TextValue "D35" "0.7598"
Set-TextValue "E35" "  -0.58%  "
Set-TextValue "D36" "1.150"
Set-TextValue "E36" "  -2.48%  "
Set-TextValue "E37" "  +0.10%  "
Set-TextValue "D38" "2.731"
Set-TextValue "E38" "  +0.22%  "
Set-TextValue "D39" "0.02013"
Set-TextValue "E39" "  -0.63%  "
Set-TextValue "D40" "80.24"
Set-TextValue "E40" "  +7.91%  "
Set-TextValue "D41" "2.662"
Set-TextValue "E41" "  -1.35%  "
Set-TextValue "D42" "6.511"
Set-TextValue "E42" "  +0.29%  "
Set-TextValue "D43" "2.085"
Set-TextValue "E43" "  -3.36%  "
Set-TextValue "D44" "0.9037"
Set-TextValue "E44" "  +2.34%  "
Set-TextValue "D45" "0.4438"
Set-TextValue "E45" "  -1.13%  "
Set-TextValue "D46" "108.10"
Set-TextValue "E46" "  -1.00%  "
Set-TextValue "D47" "1.000"
Set-TextValue "E47" "  +0.19%  "
Set-TextValue "D48" "7.806"
Set-TextValue "E48" "  -3.57%  "
Set-TextValue "D49" "1.005.40"
Set-TextValue "E49" "  +2.77%  "
Set-TextValue "D50" "0.1248"
Set-TextValue "E50" "  -1.31%  "
Set-TextValue "D51" "36.18"
Set-TextValue "E51" "  +1.28%  "
